$d = $word.ActiveDocument

for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $sec = $d.Sections.Item($si)

    for ($hi = 1; $hi -le $sec.Headers.Count; $hi++) {
        $hf = $sec.Headers.Item($hi)
        if ($hf.Exists) {
            $shapes = $hf.Range.InlineShapes
            for ($k = 1; $k -le $shapes.Count; $k++) {
                $s = $shapes.Item($k)
                $alt = $s.AlternativeText
                if ($alt -eq "BTec_Logo-Orange") {
                    $s.Name = "image2.jpg"
                } elseif ($alt -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $s.Name = "image1.png"
                }
            }
        }
    }

    for ($fi = 1; $fi -le $sec.Footers.Count; $fi++) {
        $ft = $sec.Footers.Item($fi)
        if ($ft.Exists) {
            $shapes = $ft.Range.InlineShapes
            for ($k = 1; $k -le $shapes.Count; $k++) {
                $s = $shapes.Item($k)
                $alt = $s.AlternativeText
                if ($alt -eq "BTec_Logo-Orange") {
                    $s.Name = "image2.jpg"
                } elseif ($alt -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $s.Name = "image1.png"
                }
            }
        }
    }
}

for ($k = 1; $k -le $d.InlineShapes.Count; $k++) {
    $s = $d.InlineShapes.Item($k)
    $alt = $s.AlternativeText
    if ($alt -eq "BTec_Logo-Orange") {
        $s.Name = "image2.jpg"
    } elseif ($alt -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
        $s.Name = "image1.png"
    }
}
